# Update crypto price/volume data in Sheet1 (A1:E51 table).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.463.40"
$ws.Range("E2").Value = "  -1.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.431.84"
$ws.Range("E3").Value = "  -2.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.11"
$ws.Range("E5").Value = "  -1.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.27"
$ws.Range("E6").Value = "  -3.97%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").Value = "  -2.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.427.91"
$ws.Range("E9").Value = "  -2.25%  "

$ws.Range("E10").Value = "  -4.97%  "

$ws.Range("E11").Value = "  +1.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.18"
$ws.Range("E12").Value = "  -2.92%  "

$ws.Range("E13").Value = "  -3.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.54"
$ws.Range("E14").Value = "  -2.82%  "

$ws.Range("E15").Value = "  -5.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.871.70"
$ws.Range("E16").Value = "  -2.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.417.54"
$ws.Range("E17").Value = "  -1.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.439.32"
$ws.Range("E18").Value = "  -1.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.11"
$ws.Range("E19").Value = "  -4.14%  "

$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.17"

$ws.Range("E22").Value = "  -2.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.05"
$ws.Range("E23").Value = "  +8.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -4.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.00"
$ws.Range("E25").Value = "  -3.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "611.34"
$ws.Range("E26").Value = "  -4.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.89"
$ws.Range("E27").Value = "  +1.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0973"
$ws.Range("E28").Value = "  -7.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.553.65"
$ws.Range("E29").Value = "  -3.88%  "

$ws.Range("E30").Value = "  +0.36%  "

$ws.Range("E31").Value = "  -3.99%  "

$ws.Range("E32").Value = "  -4.68%  "

$ws.Range("E33").Value = "  -2.29%  "

$ws.Range("E34").Value = "  -5.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.02"
$ws.Range("E35").Value = "  -3.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.48"
$ws.Range("E36").Value = "  -4.73%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("E38").Value = "  -3.21%  "

$ws.Range("E39").Value = "  -1.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "147.23"
$ws.Range("E40").Value = "  +0.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.23"
$ws.Range("E41").Value = "  -4.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.73"
$ws.Range("E42").Value = "  -6.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.53"
$ws.Range("E44").Value = "  -3.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.08"
$ws.Range("E45").Value = "  +0.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.65"
$ws.Range("E46").Value = "  -4.06%  "

$ws.Range("E47").Value = "  -1.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.25"
$ws.Range("E48").Value = "  -4.39%  "

$ws.Range("E49").Value = "  -4.42%  "

$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("E51").Value = "  -4.79%  "
